# Regenerate the "K" column (column G) values for rows 3-44 on the active
# worksheet. These figures come from a fresh recompute of strikeouts ("K")
# to replace the older "Strike#" derived values, per the commit message
# "regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newK = [ordered]@{
    3  = 3
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 3
    11 = 1
    12 = 0
    13 = 2
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 2
    20 = 1
    21 = 2
    22 = 2
    23 = 0
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 2
    31 = 0
    32 = 2
    33 = 1
    34 = 0
    35 = 0
    36 = 4
    37 = 2
    38 = 1
    39 = 1
    40 = 1
    41 = 4
    42 = 5
    43 = 1
    44 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
